$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    12 = 0
    13 = 1
    15 = 2
    16 = 1
    17 = 2
    18 = 1
    19 = 2
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    27 = 0
    28 = 2
    29 = 2
    30 = 0
    31 = 3
    32 = 2
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    46 = 2
    47 = 2
    48 = 2
    49 = 2
    50 = 1
    52 = 0
    53 = 1
    54 = 1
    55 = 1
    56 = 3
    57 = 2
    60 = 3
    61 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
